$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "Name_cas_test"
$ws.Range("A1").Select()
